$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("attendance")

# Mark attendance ("P") for students present in this lab session.
$presentRows = @(3, 5, 7, 8, 9, 10, 11, 12, 13, 15)
foreach ($r in $presentRows) {
    $ws.Cells.Item($r, 3).Value = "P"
}

# A few students moved between groups - update their display names.
$ws.Range("B6").Value = "Dulau I. Marius Cristian (mutat gr 4)"
$ws.Range("B18").Value = "Prata L. Dragos Liviu (mutat gr 4)"
$ws.Range("B24").Value = "Vranau V. Flavius Silviu (mutat gr 4)"

# A student moved in from group 4 - add him as a new row.
$ws.Range("A26").Value = 25
$ws.Range("B26").Value = "Bompa Flaviu (venit din gr 4)"
$ws.Range("C26").Value = "P"

# Update the active selection to reflect where edits were made.
$ws.Range("B6").Select()
